$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("transactions")
$ws.Activate()

# Clear the placeholder "0/0/0000" dates (literal value in D101/E101, and the
# formulas that propagated it down through D102:D150 / E102:E150) while
# keeping the existing cell formatting (style "8") intact.
$ws.Range("D101:E150").ClearContents()

# Restore the view: scroll the window so row 84 is at the top and move the
# active selection to D101:E150 (anchored at D101).
$excel.ActiveWindow.ScrollRow = 84
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("D101:E150").Select()
$excel.ActiveCell = $ws.Range("D101")
